# IfThenBlocks.pptx update:
#   - refresh the "last edited" auto-date fields (handout + notes masters)
#   - rename slide 2 title: "Lesson Objectives  UPDATE!" -> "Lesson Objectives"
#   - rename slide 3 title: "Repeating code" -> "If Then Blocks"

$p = $ppt.ActivePresentation

# --- Handout master date placeholder ---
$handoutMaster = $p.HandoutMaster
$handoutDate = $handoutMaster.HeadersFooters.DateAndTime
$handoutDate.Text = "11/23/2020"

# --- Notes master date placeholder ---
$notesMaster = $p.NotesMaster
$notesDate = $notesMaster.HeadersFooters.DateAndTime
$notesDate.Text = "11/23/2020"

# --- Slide 2 title ---
$slide2 = $p.Slides.Item(2)
$slide2.Shapes.Item(1).TextFrame.TextRange.Text = "Lesson Objectives"

# --- Slide 3 title ---
$slide3 = $p.Slides.Item(3)
$slide3.Shapes.Item(1).TextFrame.TextRange.Text = "If Then Blocks"
